$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Day 4" row of data (row 5)
$ws.Range("A5").Value = "Day 4"
$ws.Range("B5").Value = 45806
$ws.Range("C5").Value = "Binary Search"
$ws.Range("D5").Value = "Best Time to Buy and Sell Stock"
$ws.Range("E5").Value = "Longest Substring Without Repeating Characters"
$ws.Range("F5").Value = "Binary Search, Greedy, Sliding Window, Hash Set"
$ws.Range("G5").Value = "S"
$ws.Range("H5").Value = "YES"

# Match the date number formatting used by the rows above (reuse existing style)
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Widen column E to fit the new, longer topic text
$ws.Columns.Item(5).ColumnWidth = 36.67

# Update the active selection to D6, as recorded in the saved workbook
$null = $ws.Range("D6").Select()
